# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) used the literal text "4-25-2012-13" for every
# team row; correct it to the ISO-ish form "2013-04-25".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # Column BF = 58
    if ($cell.Value2 -eq "4-25-2012-13") {
        # Leading apostrophe forces text so Excel doesn't reinterpret the
        # ISO-looking "2013-04-25" as a date serial number; the stored
        # cell value/text is the plain string without the apostrophe.
        $cell.Value = "'2013-04-25"
    }
}
